# Hortaliza, Macroferia Regional de Talca - Coliflor: weekly refresh.
# Two new price observations are inserted at the top of the data block
# (rows 251-252), pushing every existing record down by two rows
# (old A1:R279 -> new A1:R281).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 251, shifting the
# rest of the table (old rows 251..279) down to 253..281.
$ws.Range("A251:R252").EntireRow.Insert()

# --- New row 251 -----------------------------------------------------
$ws.Range("A251").Value = 5
$ws.Range("B251").Value = "Macroferia Regional de Talca"
$ws.Range("C251").Value = "Maule"
$ws.Range("D251").Value = 44769
$ws.Range("E251").Value = 7
$ws.Range("F251").Value = 100112008
$ws.Range("G251").Value = "Coliflor"
$ws.Range("H251").Value = "Sin especificar"
$ws.Range("I251").Value = "Primera"
$ws.Range("J251").Value = 3000
$ws.Range("K251").Value = 1000
$ws.Range("L251").Value = 1000
$ws.Range("M251").Value = 1000
$ws.Range("N251").Value = "`$/unidad"
$ws.Range("O251").Value = "Región del Maule"
$ws.Range("P251").Value = 1000
$ws.Range("Q251").Value = 1
$ws.Range("R251").Value = "Hortaliza"

# --- New row 252 -----------------------------------------------------
$ws.Range("A252").Value = 5
$ws.Range("B252").Value = "Macroferia Regional de Talca"
$ws.Range("C252").Value = "Maule"
$ws.Range("D252").Value = 44769
$ws.Range("E252").Value = 7
$ws.Range("F252").Value = 100112008
$ws.Range("G252").Value = "Coliflor"
$ws.Range("H252").Value = "Sin especificar"
$ws.Range("I252").Value = "Segunda"
$ws.Range("J252").Value = 3000
$ws.Range("K252").Value = 800
$ws.Range("L252").Value = 800
$ws.Range("M252").Value = 800
$ws.Range("N252").Value = "`$/unidad"
$ws.Range("O252").Value = "Región del Maule"
$ws.Range("P252").Value = 800
$ws.Range("Q252").Value = 1
$ws.Range("R252").Value = "Hortaliza"

# Match the date-cell formatting used by the rest of column D (style
# index 2 in the original workbook -> numFmt "YYYY-MM-DD HH:MM:SS").
$ws.Range("D251:D252").NumberFormat = $ws.Range("D253").NumberFormat
